$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (shifts existing rows 32-48 down to 33-49)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly price record
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44553
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112021
$ws.Range("G32").Value = "Ají"
$ws.Range("H32").Value = "Americana (o)"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = 15250
$ws.Range("N32").Value = "`$/caja 14 kilos"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 1089
$ws.Range("Q32").Value = 14
$ws.Range("R32").Value = "Hortaliza"
